$d = $word.ActiveDocument

# The edit removes three paragraphs that followed the
# "LOQ4047: Trabalho de Conclusão de Curso I (Requisito)" paragraph:
#   1) a blank paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#       pages. Original theme under Creative Commons Attribution"
# leaving the blank paragraph that precedes the page-break paragraph intact.

$marker = "LOQ4047: Trabalho de Conclusão de Curso I (Requisito)"

$startPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*$marker*") {
        $startPara = $p
        break
    }
}

$removeFrom = $startPara.Next()
$removeTo = $removeFrom.Next().Next()

$d.Range($removeFrom.Range.Start, $removeTo.Range.End).Delete()
